# oh my god big work
# Remove the hard-coded "2021年" placeholder text in three spots so the
# template only keeps the bracketed fill-in placeholders.
#
#   1) "...开发中心于2021年【起始日期】至..."      -> "...开发中心于【起始日期】至..."
#   2) "...报告日期：2021年【出报告日期】）。"        -> "...报告日期：【出报告日期】）。"
#   3) "上海计算机软件技术开发中心\n2021年【出报告日期】" -> "...\n【出报告日期】"

$d = $word.ActiveDocument

function Remove-LiteralText($context, $target) {
    # Find a unique chunk of text ($context) that contains $target, then
    # delete just the $target sub-range inside it (by character offset),
    # leaving the runs before/after $context's match completely untouched.
    $rng = $d.Content
    $found = $rng.Find.Execute($context, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        Write-Output "NOT FOUND: $context"
        return $false
    }
    $offset = $context.IndexOf($target)
    if ($offset -lt 0) {
        Write-Output "TARGET NOT IN CONTEXT: $target"
        return $false
    }
    $start = $rng.Start + $offset
    $end = $start + $target.Length
    $victim = $d.Range($start, $end)
    $victim.Delete()
    return $true
}

function Remove-LiteralTextNoMerge($context, $target) {
    # Same idea as Remove-LiteralText, but used when deleting $target would
    # leave two neighboring runs that share identical run formatting (rPr) —
    # Word's delete-merge logic would otherwise splice those two runs into
    # one. Toggling Bold off->on->off on the trailing run (a no-op on the
    # rendered formatting) forces it to stay a distinct run.
    $rng = $d.Content
    $found = $rng.Find.Execute($context, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        Write-Output "NOT FOUND: $context"
        return $false
    }
    $offset = $context.IndexOf($target)
    if ($offset -lt 0) {
        Write-Output "TARGET NOT IN CONTEXT: $target"
        return $false
    }
    $base = $rng.Start
    $start = $base + $offset
    $end = $start + $target.Length

    $tailStart = $end
    $tailEnd = $base + $context.Length
    $tailLen = $tailEnd - $tailStart

    if ($tailLen -gt 0) {
        $tailRng = $d.Range($tailStart, $tailEnd)
        $origBold = $tailRng.Bold
        $tailRng.Bold = 1
    }

    $victim = $d.Range($start, $end)
    $victim.Delete()

    if ($tailLen -gt 0) {
        $newTailRng = $d.Range($start, $start + $tailLen)
        $newTailRng.Bold = $origBold
    }
    return $true
}

Remove-LiteralText "于2021年【起始日期】" "2021年" | Out-Null
Remove-LiteralTextNoMerge "报告日期：2021年【出报告日期】" "2021年" | Out-Null
Remove-LiteralText "2021年【出报告日期】" "2021年" | Out-Null

Write-Output "done"
